$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3715635136444516
$ws.Range("C2").Value = 0.1535251409344947
$ws.Range("D2").Value = -0.3787274915688817

$ws.Range("B3").Value = 0.3978800233787532
$ws.Range("C3").Value = -0.5096886910384966

$ws.Range("B4").Value = -0.4316193217044243

$ws.Range("B5").Value = -0.1122912858586644
$ws.Range("C5").Value = 0.3447545960399119
$ws.Range("D5").Value = 0.227751383870373
$ws.Range("E5").Value = 0.2766951351148101

$ws.Range("B6").Value = 0.2910639859155906
$ws.Range("C6").Value = 0.1632992805275675
$ws.Range("D6").Value = 0.2326774114107624

$ws.Range("B7").Value = 0.2439681181820851
$ws.Range("C7").Value = 0.2460476392559225

$ws.Range("B8").Value = 0.1769863131018627

$ws.Range("B9").Value = -0.1631350120513026
$ws.Range("C9").Value = 0.1559839361989813
$ws.Range("D9").Value = 0.2182565853642597
$ws.Range("E9").Value = -0.258310573012756

$ws.Range("B10").Value = 0.1101846028489047
$ws.Range("C10").Value = 0.233606398273979
$ws.Range("D10").Value = -0.2104585808681099

$ws.Range("B11").Value = 0.5197732891916149
$ws.Range("C11").Value = -0.2505245903104478

$ws.Range("B12").Value = -0.2267575370288383

$ws.Range("B13").Value = -0.1046767183175807
$ws.Range("C13").Value = -0.08657733350508284
$ws.Range("D13").Value = -0.2473038116165672

$ws.Range("B14").Value = 0.04486822346272668
$ws.Range("C14").Value = -0.290606343338754

$ws.Range("B15").Value = -0.5174346182871943
